$d = $word.ActiveDocument

# 1. Merge the split runs "...interface for grounds that" + "are able to" + "be jumped"
#    (with a gramStart/gramEnd proofErr pair in between) into a single run/text node.
#    Word's Find/Replace merges the runs spanned by the match into one run.
$find1 = $d.Content.Find
$find1.Execute(" interface for grounds that are able to be jumped", $false, $false, $false, $false, $false, $true, 1, $false, " interface for grounds that are able to be jumped", 2)

# 2. Remove the stray empty paragraph between "Signed by (type..." and
#    "I, Ong Di Sheng accept this WBA".
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -match "^Signed by") {
        $blank = $paras.Item($i + 1)
        if ($blank.Range.Text.Trim() -eq "") {
            $blank.Range.Delete()
        }
        break
    }
}

# 3. Drop the stale lastRenderedPageBreak cached on Kennedy Tan Sing Ye's
#    acceptance paragraph (a re-find/replace over the text normalizes the
#    run and clears the rendering cache marker).
$find2 = $d.Content.Find
$find2.Execute("I, Kennedy Tan Sing Ye accept this WBA", $false, $false, $false, $false, $false, $true, 1, $false, "I, Kennedy Tan Sing Ye accept this WBA", 2)
